# Time Log for CISC 4900 - add new "1029 - 1112" sheet
# -------------------------------------------------------------------
# Strategy: duplicate the most recent weekly sheet ("1015 - 1028") so
# that the new sheet inherits the same column widths, fonts, borders,
# number formats and overall look & feel, then rename it, strip the
# bits that shouldn't carry over (the old hyperlink), overwrite the
# log rows with the new entries, extend it with one more blank
# placeholder row, rebuild the data validation + table over the new
# range and fix up the "Total Duration" formula to point at the new
# table.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("1015 - 1028")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$template.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "1029 - 1112"

# Drop the hyperlink that was copied over from the template sheet.
foreach ($h in $ws.Hyperlinks) {
    $h.Delete()
}

# ---- Overwrite the log entries -------------------------------------------------
$ws.Range("A2").Value2 = 45599
$ws.Range("B2").Value2 = 1
$ws.Range("C2").Value2 = "Research"
$ws.Range("D2").Value2 = "looked for ways to contain Alaska and Hawaii in seperate containers and still be able to use them as filters"
$ws.Range("E2").Value2 = ""

$ws.Range("A3").Value2 = 45599
$ws.Range("B3").Value2 = 6
$ws.Range("C3").Value2 = "Dashboard"
$ws.Range("D3").Value2 = "Built upon the Covid 2020 Death Report. It now includes Infection as well. Modernized the design of that dashboard as well."
$ws.Range("E3").Value2 = ""

$ws.Range("A4").Value2 = 45600
$ws.Range("B4").Value2 = 2
$ws.Range("C4").Value2 = "Research"
$ws.Range("D4").Value2 = "Looked for ways to implement Time Series Map Chart. Committed to repo."
$ws.Range("E4").Value2 = ""

$ws.Range("A5").Value2 = 45602
$ws.Range("B5").Value2 = 2
$ws.Range("C5").Value2 = "Research"
$ws.Range("D5").Value2 = "looked for ways to deploy tableau. Currently most reasonable option is through Tableau public and embed unto a website."
$ws.Range("E5").Value2 = ""

$ws.Range("A6").Value2 = 45603
$ws.Range("B6").Value2 = 8
$ws.Range("C6").Value2 = "Dashboard"
$ws.Range("D6").Value2 = "Adjusted misinterpreted data on Covid Death/Infection Report. Started on a Time Series Map Chart Dashboard."
$ws.Range("E6").Value2 = ""

$ws.Range("A7").Value2 = 45605
$ws.Range("B7").Value2 = 7
$ws.Range("C7").Value2 = "Dashboard"
$ws.Range("D7").Value2 = "Finished Time Series Dashboard. Started on Reworking Political Affiliations Dashboard"
$ws.Range("E7").Value2 = ""

$ws.Range("A8").Value2 = 45606
$ws.Range("B8").Value2 = 0.5
$ws.Range("C8").Value2 = "Jira"
$ws.Range("D8").Value2 = "Updates to Jira Project Management Board"
$ws.Range("E8").Value2 = ""

$ws.Range("A9").Value2 = 45606
$ws.Range("B9").Value2 = 1
$ws.Range("C9").Value2 = "Demo Recording"
$ws.Range("D9").Value2 = "created a baseline script to use to record 2nd demo"
$ws.Range("E9").Value2 = "Demo Recorded"

$ws.Range("A10").Value2 = 45606
$ws.Range("B10").Value2 = 0.5
$ws.Range("C10").Value2 = "Email"
$ws.Range("D10").Value2 = "Emailed Supervisor regarding missing interim evaulation report"
$ws.Range("E10").Value2 = ""

$ws.Range("A11").Value2 = 45606
$ws.Range("B11").Value2 = 2
$ws.Range("C11").Value2 = "Dashboard"
$ws.Range("D11").Value2 = "Took some inspirations from other dashboards on Tableau Public for data visualization. Changed overall color theme of dashboards to black to help dashboard visualization. Adjusted dashboard heights, widiths, and inner/outer paddings to ensure visualization consistency."
$ws.Range("E11").Value2 = ""

# ---- Add one extra blank placeholder row (rows 12-15 already came over blank
# from the template; append row 16 using the same formatting as row 15). ------
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16:E16").ClearContents()
$excel.CutCopyMode = 0

# ---- Data validation now spans the extra row ------------------------------
$dvRange = $ws.Range("A2:A16")
$dvRange.Validation.Delete()
$dvRange.Validation.Add(7, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '=OR(NOT(ISERROR(DATEVALUE(A2))), AND(ISNUMBER(A2), LEFT(CELL("format", A2))="D"))')
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $false
$dvRange.Validation.ShowInput = $false
$dvRange.Validation.ShowError = $false

# ---- Rebuild the table over the new sheet ---------------------------------
$tableRange = $ws.Range("A1:E16")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table2_2"
$lo.TableStyle = "1029 - 1112-style"
$lo.ShowTableStyleFirstColumn = $true
$lo.ShowTableStyleLastColumn = $true
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowAutoFilter = $false
$lo.ShowTotals = $false

# ---- Fix the "Total Duration" formula to reference the new table ----------
$ws.Range("H3").Formula = "=SUM(Table2_2[Duration (hours)])"

$ws.Activate()
